$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room: the old blank spacer row 139 needs to become a real time-tracking
# row, and two *additional* blank spacer rows are inserted after it (so the
# totals block that used to start at row 140 now starts at row 143).
$ws.Range("139:141").Insert()

# Rows 140 and 141 are brand-new blank spacer rows - keep only the D/F/G
# placeholders (same look as the other spacer rows), drop everything else
# that got copied down from row 139 by the insert.
$ws.Range("A140:C141").Clear()
$ws.Range("E140:E141").Clear()
$ws.Range("H140:K141").Clear()

# Row 139: same "Interface Design / MockUps / [FEATURE]" entry pattern as the
# rows above it, now for the new "Einstellungsicons" task.
$ws.Cells.Item(139,1).Value = 22
$ws.Cells.Item(139,2).Value = "Interface Design"
$ws.Cells.Item(139,3).Value = "MockUps"
$ws.Cells.Item(139,4).Value = "[FEATURE]"
$ws.Cells.Item(139,4).Font.Color = 0
$ws.Cells.Item(139,5).Value = "Einstellungsicons"
$ws.Cells.Item(139,6).Value = 44464
$ws.Cells.Item(139,6).Font.Color = 0
$ws.Cells.Item(139,7).Value = 44481
$ws.Cells.Item(139,7).Font.Color = 0
$ws.Cells.Item(139,8).Font.Color = 0

$ws.Cells.Item(139,9).Formula = "=ROUNDUP(((SUM(K139-J139)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(139,10).Formula = "=K138"

# New "K" column look for this entry: light grey fill, black centred time text.
$k = $ws.Cells.Item(139,11)
$k.Value = 0.53125
$k.NumberFormat = "h:mm"
$k.Font.Color = 0
$k.Interior.Color = 15921906
$k.HorizontalAlignment = -4108

$ws.Range("H138").Select()
